# Apply this edit using the Excel COM object model.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (numeric-looking price/volume strings) to be stored
# as text so Excel does not silently coerce values like "1.013" into numbers.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = "27.618.90"
$ws.Cells.Item(2, 5).Value = "  -0.91%  "

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = "1.858.74"
$ws.Cells.Item(3, 5).Value = "  -1.50%  "

# Row 4: TetherUSD
$ws.Cells.Item(4, 4).Value = "1.013"
$ws.Cells.Item(4, 5).Value = "  -0.37%  "

# Row 5: BNB
$ws.Cells.Item(5, 4).Value = "334.97"
$ws.Cells.Item(5, 5).Value = "  -0.18%  "

# Row 6: USDC
$ws.Cells.Item(6, 4).Value = "1.012"
$ws.Cells.Item(6, 5).Value = "  -0.41%  "

# Row 7: XRP
$ws.Cells.Item(7, 4).Value = "0.4637"
$ws.Cells.Item(7, 5).Value = "  -1.03%  "

# Row 8: Cardano
$ws.Cells.Item(8, 4).Value = "0.3916"
$ws.Cells.Item(8, 5).Value = "  +0.05%  "

# Row 9: OKB
$ws.Cells.Item(9, 4).Value = "46.04"
$ws.Cells.Item(9, 5).Value = "  -3.41%  "

# Row 10: Dogecoin
$ws.Cells.Item(10, 4).Value = "0.07963"
$ws.Cells.Item(10, 5).Value = "  -1.03%  "

# Row 11: Polygon
$ws.Cells.Item(11, 4).Value = "0.9988"
$ws.Cells.Item(11, 5).Value = "  -1.93%  "

# Row 12: Solana
$ws.Cells.Item(12, 4).Value = "21.61"
$ws.Cells.Item(12, 5).Value = "  -0.86%  "

# Row 13: WrappedEther
$ws.Cells.Item(13, 4).Value = "1.863.12"
$ws.Cells.Item(13, 5).Value = "  -1.23%  "

# Row 14: Polkadot
$ws.Cells.Item(14, 4).Value = "5.942"
$ws.Cells.Item(14, 5).Value = "  -0.25%  "

# Row 15: Chainlink
$ws.Cells.Item(15, 4).Value = "7.202"
$ws.Cells.Item(15, 5).Value = "  +1.53%  "

# Row 16: BinanceUSD
$ws.Cells.Item(16, 4).Value = "1.014"
$ws.Cells.Item(16, 5).Value = "  -0.24%  "

# Row 17: Litecoin
$ws.Cells.Item(17, 4).Value = "88.24"
$ws.Cells.Item(17, 5).Value = "  +1.22%  "

# Row 18: TRON
$ws.Cells.Item(18, 4).Value = "0.06712"
$ws.Cells.Item(18, 5).Value = "  -0.83%  "

# Row 19: ShibaInu
$ws.Cells.Item(19, 4).Value = "0.00001040"
$ws.Cells.Item(19, 5).Value = "  -0.73%  "

# Row 20: Avalanche
$ws.Cells.Item(20, 5).Value = "  -0.04%  "

# Row 21: Dai
$ws.Cells.Item(21, 4).Value = "1.012"
$ws.Cells.Item(21, 5).Value = "  -0.30%  "

# Row 22: WrappedBTC
$ws.Cells.Item(22, 4).Value = "27.614.27"
$ws.Cells.Item(22, 5).Value = "  -0.86%  "

# Row 23: Uniswap
$ws.Cells.Item(23, 4).Value = "5.444"
$ws.Cells.Item(23, 5).Value = "  -1.00%  "

# Row 24: Cosmos
$ws.Cells.Item(24, 4).Value = "10.94"
$ws.Cells.Item(24, 5).Value = "  -0.50%  "

# Row 25: Toncoin
$ws.Cells.Item(25, 4).Value = "2.302"
$ws.Cells.Item(25, 5).Value = "  -1.86%  "

# Row 26: WrappedliquidstakedEther2.0
$ws.Cells.Item(26, 4).Value = "2.078.46"
$ws.Cells.Item(26, 5).Value = "  -1.45%  "

# Row 27: Monero
$ws.Cells.Item(27, 4).Value = "159.61"
$ws.Cells.Item(27, 5).Value = "  -0.03%  "

# Row 28: EthereumClassic
$ws.Cells.Item(28, 4).Value = "19.63"
$ws.Cells.Item(28, 5).Value = "  -2.09%  "

# Row 29: LidoDAOToken
$ws.Cells.Item(29, 4).Value = "2.138"
$ws.Cells.Item(29, 5).Value = "  +3.24%  "

# Row 30: InternetComputer(DFINITY)
$ws.Cells.Item(30, 4).Value = "5.419"
$ws.Cells.Item(30, 5).Value = "  -0.84%  "

# Row 31: BitcoinCash
$ws.Cells.Item(31, 4).Value = "121.85"
$ws.Cells.Item(31, 5).Value = "  -0.02%  "

# Row 32: ImmutableX
$ws.Cells.Item(32, 4).Value = "0.9754"
$ws.Cells.Item(32, 5).Value = "  +1.12%  "

# Row 33: Stellar
$ws.Cells.Item(33, 4).Value = "0.09420"
$ws.Cells.Item(33, 5).Value = "  -0.70%  "

# Row 34: HuobiToken
$ws.Cells.Item(34, 4).Value = "3.626"
$ws.Cells.Item(34, 5).Value = "  -0.64%  "

# Row 35: Filecoin
$ws.Cells.Item(35, 4).Value = "5.305"
$ws.Cells.Item(35, 5).Value = "  -0.61%  "

# Row 36: ARBITRUM
$ws.Cells.Item(36, 5).Value = "  -4.79%  "

# Row 37: VeChain
$ws.Cells.Item(37, 5).Value = "  -0.67%  "

# Row 38: Hedera
$ws.Cells.Item(38, 4).Value = "0.06012"
$ws.Cells.Item(38, 5).Value = "  -1.51%  "

# Row 39: FraxShare
$ws.Cells.Item(39, 4).Value = "8.316"
$ws.Cells.Item(39, 5).Value = "  +3.49%  "

# Row 40: TrustWalletToken
$ws.Cells.Item(40, 4).Value = "1.195"
$ws.Cells.Item(40, 5).Value = "  -1.65%  "

# Row 41: TheSandbox -> Frax
$ws.Cells.Item(41, 2).Value = "Frax"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(41, 4).Value = "1.012"
$ws.Cells.Item(41, 5).Value = "  -0.29%  "

# Row 42: Algorand -> TheSandbox
$ws.Cells.Item(42, 2).Value = "TheSandbox"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(42, 4).Value = "0.5936"
$ws.Cells.Item(42, 5).Value = "  -0.57%  "

# Row 43: Aptos -> Algorand
$ws.Cells.Item(43, 2).Value = "Algorand"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(43, 4).Value = "0.1866"
$ws.Cells.Item(43, 5).Value = "  -0.83%  "

# Row 44: WEMIXTOKEN -> Aptos
$ws.Cells.Item(44, 2).Value = "Aptos"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(44, 4).Value = "10.34"
$ws.Cells.Item(44, 5).Value = "  +0.80%  "

# Row 45: Decentraland -> WEMIXTOKEN
$ws.Cells.Item(45, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(45, 4).Value = "1.248"
$ws.Cells.Item(45, 5).Value = "  -1.59%  "

# Row 46: EnergySwap -> Decentraland
$ws.Cells.Item(46, 2).Value = "Decentraland"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(46, 4).Value = "0.5598"
$ws.Cells.Item(46, 5).Value = "  -1.65%  "

# Row 47: NEARProtocol -> EnergySwap
$ws.Cells.Item(47, 2).Value = "EnergySwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47, 4).Value = "12.08"
$ws.Cells.Item(47, 5).Value = "  -0.10%  "

# Row 48: Cronos -> NEARProtocol
$ws.Cells.Item(48, 2).Value = "NEARProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(48, 4).Value = "1.915"
$ws.Cells.Item(48, 5).Value = "  -0.64%  "

# Row 49: Quant -> Cronos
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(49, 4).Value = "0.06714"
$ws.Cells.Item(49, 5).Value = "  -3.03%  "

# Row 50: EOS -> Quant
$ws.Cells.Item(50, 2).Value = "Quant"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(50, 4).Value = "111.60"
$ws.Cells.Item(50, 5).Value = "  -1.93%  "

# Row 51: PaxDollar -> EOS
$ws.Cells.Item(51, 2).Value = "EOS"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Cells.Item(51, 4).Value = "1.050"
$ws.Cells.Item(51, 5).Value = "  -1.62%  "

# Restore the default (unstyled) cell style now that the values are text,
# matching the original workbook formatting for these cells.
$priceRange.Style = "Normal"